$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Rerun the example data": the `group` column (J2:J17) was regenerated and its
# value changed from "2" to "3" for every data row. The column holds text
# (shared-string) values, so re-enter it with a leading apostrophe to keep it
# stored as text rather than being auto-coerced to a number.
$ws.Range("J2:J17").Value = "'3"

# The Collectiondate column (C2:C17) was previously formatted with a custom
# number format (mm/dd/yyyy, numFmtId 165). After the rerun it uses Excel's
# built-in short-date format instead of the custom one.
$ws.Range("C2:C17").NumberFormat = "mm-dd-yy"
